# Adjust Investment Summary (and Timeline & Milestones) table column widths.
#
# PowerPoint COM table Column.Width is expressed in points; OOXML a:gridCol
# widths are stored in EMUs (1 pt = 12700 EMU). Setting every column's width
# causes PowerPoint to recompute the graphicFrame's a:ext cx as the sum of
# the grid column widths, matching the target diff exactly.

$p = $ppt.ActivePresentation

# --- Slide 5: "Timeline & Milestones" table (4 columns) ---
$s1 = $p.Slides.Item(5)
$t1 = $s1.Shapes.Item(3).Table

$t1.Columns.Item(1).Width = 871093 / 12700
$t1.Columns.Item(2).Width = 2177733 / 12700
$t1.Columns.Item(3).Width = 1306639 / 12700
$t1.Columns.Item(4).Width = 4355466 / 12700

# --- Slide 8: "Investment Summary" table (7 columns) ---
$s2 = $p.Slides.Item(8)
$t2 = $s2.Shapes.Item(3).Table

$t2.Columns.Item(1).Width = 1742186 / 12700
$t2.Columns.Item(2).Width = 1045311 / 12700
$t2.Columns.Item(3).Width = 2003514 / 12700
$t2.Columns.Item(4).Width = 1132421 / 12700
$t2.Columns.Item(5).Width = 871093 / 12700
$t2.Columns.Item(6).Width = 871093 / 12700
$t2.Columns.Item(7).Width = 1045311 / 12700
